# Append the 2025-03-24 price row (row 23) to every price sheet in the
# Solar_Prices workbook. Each sheet gets the same new date in column A and
# a sheet-specific price in column B, matching the previous day's format.
#
# The new date/price are entered as literal text (not auto-converted to a
# date serial / number) so the stored values match the existing rows,
# which are themselves plain text - hence the NumberFormat "@" (Text)
# applied before assigning the value.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-24"

# Ordered (sheet name, new price) pairs - one per price sheet in the workbook.
$sheetPrices = @(
    ,("N-Dense",                 "40")
    ,("N-Type",                  "43")
    ,("N-type Wafer",            "1.19")
    ,("Cell Topcon 183mm",       "0.298")
    ,("Module Topcon 183mm",     "0.1")
    ,("Silver Rear_side",        "5,399")
    ,("Silver Busbar front-side","8,083")
    ,("Silver finger front-side","8,133")
    ,("USD_CNY",                 "7.2717")
)

foreach ($pair in $sheetPrices) {
    $sheetName = $pair[0]
    $price = $pair[1]

    $ws = $wb.Worksheets.Item($sheetName)

    $dateCell = $ws.Range("A23")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate

    $priceCell = $ws.Range("B23")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
}
